$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 189, pushing existing rows 189-212 down to 190-213.
$ws.Rows.Item(189).Insert()

# Populate the new row 189 with the new record's data (same fixed
# columns as the rest of the table, new price/date data).
$ws.Range("A189").Value = 4
$ws.Range("B189").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C189").Value = "Los Lagos"
$ws.Range("D189").Value = 44694
$ws.Range("E189").Value = 10
$ws.Range("F189").Value = 100112039
$ws.Range("G189").Value = "Ciboulette"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 240
$ws.Range("K189").Value = 2500
$ws.Range("L189").Value = 2500
$ws.Range("M189").Value = 2500
$ws.Range("N189").Value = "$/docena de atados"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 833
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"
